$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.05873514465104036
$ws.Range("C2").Value = 0.0009034841054524714
$ws.Range("B3").Value = 0.1637156252749384
$ws.Range("C3").Value = 0.002126261199068803
